$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update two existing footprint values (shared strings edited in place) ---
$ws.Range("B16").Value = "conservify:SHT3x"
$ws.Range("B92").Value = "conservify:ATLAS_SENSOR_BOARD_COMBO"
$ws.Range("B93").Value = "conservify:ATLAS_SENSOR_BOARD_COMBO"
$ws.Range("B94").Value = "conservify:ATLAS_SENSOR_BOARD_COMBO"
$ws.Range("B95").Value = "conservify:ATLAS_SENSOR_BOARD_COMBO"
$ws.Range("B96").Value = "conservify:ATLAS_SENSOR_BOARD_COMBO"

# --- Append new parts (rows 110-116) ---

# Row 110: F109, inductor (BLM18KG221SN1D) using existing IND-0603 footprint
$ws.Range("A110").Value = "F109"
$ws.Range("B110").Value = "RocketScreamKicadLibrary:IND-0603"
$ws.Range("C110").Value = "BLM18KG221SN1D"
$ws.Range("D110").Value = "81-BLM18KG221SN1D"
$ws.Range("E110").Value = "mouser"
$ws.Range("P110").Formula = '="F" & (ROW()-1)'

# Row 111: F110, inductor (BLM18KG221SN1D) with new 0603 hand-soldering footprint
$ws.Range("A111").Value = "F110"
$ws.Range("B111").Value = "Inductors_SMD:L_0603_HandSoldering"
$ws.Range("C111").Value = "BLM18KG221SN1D"
$ws.Range("D111").Value = "81-BLM18KG221SN1D"
$ws.Range("E111").Value = "mouser"

# Row 112: F111, 1x07 socket strip
$ws.Range("A112").Value = "F111"
$ws.Range("B112").Value = "Socket_Strips:Socket_Strip_Straight_1x07_Pitch2.54mm"
$ws.Range("C112").Value = "Conn_01x07"

# Row 113: F112, solder jumper (open)
$ws.Range("A113").Value = "F112"
$ws.Range("B113").Value = "conservify:SJ_OPEN"
$ws.Range("C113").Value = "Conn_01x02"
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0

# Row 114: F113, measurement point
$ws.Range("A114").Value = "F113"
$ws.Range("B114").Value = "Measurement_Points:Measurement_Point_Round-SMD-Pad_Small"
$ws.Range("C114").Value = "Conn_01x01"
$ws.Range("H114").Value = 0
$ws.Range("I114").Value = 0
$ws.Range("J114").Value = 0

# Row 115: F114, 1x03 socket strip (re-uses existing footprint text)
$ws.Range("A115").Value = "F114"
$ws.Range("B115").Value = "Socket_Strips:Socket_Strip_Straight_1x03_Pitch2.54mm"
$ws.Range("C115").Value = "Conn_01x03"
$ws.Range("H115").Value = 0
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = 0

# Row 116: F115, solder jumper (shorted / 0R)
$ws.Range("A116").Value = "F115"
$ws.Range("B116").Value = "conservify:SJ_SHORTED"
$ws.Range("C116").Value = "0R"

# --- View state: move to where the user was last working ---
$null = $ws.Range("B96").Select()
